$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Update existing values on row 2 (gold/filmgold row)
$ws.Range("B2").Value = 160
$ws.Range("C2").Value = 0.000064766799999999998

# Add a new row of data: seed
$ws.Range("A6").Value = "seed"
$ws.Range("B6").Value = 9.1999999999999993
$ws.Range("C6").Value = 0.0000034000000000000001
$ws.Range("C6").NumberFormat = $ws.Range("C2").NumberFormat

# Update selection to reflect new active cell after edits
$ws.Range("E6").Select()
